$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6255800
$ws.Range("I137").Value = 11118833
$ws.Range("J137").Value = 3328.5715
$ws.Range("K137").Value = 33356499
$ws.Range("L137").Value = 9985.7145
$ws.Range("M137").Value = -33353949
$ws.Range("N137").Value = -15085.7145
$ws.Range("H138").Value = 2219.5679
$ws.Range("I138").Value = 1311.0204
$ws.Range("J138").Value = 3610.7812
$ws.Range("K138").Value = 3933.0612
$ws.Range("L138").Value = 10832.3436
$ws.Range("M138").Value = 1206.9388
$ws.Range("N138").Value = -21112.3436
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2839.0454
$ws.Range("I61").Value = 1608.9524
$ws.Range("K61").Value = 1608.9524
$ws.Range("M61").Value = -1396.9524
$ws.Range("H74").Value = 559.5263
$ws.Range("I74").Value = 559.5263
$ws.Range("K74").Value = 559.5263
$ws.Range("M74").Value = 314.4737
$ws.Range("H77").Value = 559.5263
$ws.Range("I77").Value = 559.5263
$ws.Range("K77").Value = 2797.6315
$ws.Range("M77").Value = 1570.3685
$ws.Range("H102").Value = 3842.2222
$ws.Range("I102").Value = 3296.6667
$ws.Range("J102").Value = 4933.3335
$ws.Range("K102").Value = 3296.6667
$ws.Range("L102").Value = 4933.3335
$ws.Range("M102").Value = -1674.6667
$ws.Range("N102").Value = -8177.3335
$ws.Range("H132").Value = 15626865
$ws.Range("I132").Value = 18869296
$ws.Range("J132").Value = 4238.727
$ws.Range("K132").Value = 56607888
$ws.Range("L132").Value = 12716.181
$ws.Range("M132").Value = -56605358
$ws.Range("N132").Value = -17776.181
$ws.Range("H136").Value = 2839.0454
$ws.Range("I136").Value = 1608.9524
$ws.Range("K136").Value = 4826.857199999999
$ws.Range("M136").Value = -2276.857199999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2240.611
$ws.Range("I134").Value = 1994.6666
$ws.Range("J134").Value = 2732.5
$ws.Range("K134").Value = 5983.9998
$ws.Range("L134").Value = 8197.5
$ws.Range("M134").Value = -3448.9998
$ws.Range("N134").Value = -13267.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2568759
$ws.Range("I31").Value = 3337040
$ws.Range("J31").Value = 7822.222
$ws.Range("K31").Value = 3337040
$ws.Range("L31").Value = 7822.222
$ws.Range("M31").Value = -3336745
$ws.Range("N31").Value = -8412.222
$ws.Range("H34").Value = 2568759
$ws.Range("I34").Value = 3337040
$ws.Range("J34").Value = 7822.222
$ws.Range("K34").Value = 3337040
$ws.Range("L34").Value = 7822.222
$ws.Range("M34").Value = -3336838
$ws.Range("N34").Value = -8226.222
$ws.Range("H58").Value = 13515759
$ws.Range("I58").Value = 1224
$ws.Range("J58").Value = 41671040
$ws.Range("K58").Value = 1224
$ws.Range("L58").Value = 41671040
$ws.Range("M58").Value = -1021
$ws.Range("N58").Value = -41671446
$ws.Range("H99").Value = 2416.6667
$ws.Range("I99").Value = 625
$ws.Range("J99").Value = 3312.5
$ws.Range("K99").Value = 625
$ws.Range("L99").Value = 3312.5
$ws.Range("M99").Value = 873
$ws.Range("N99").Value = -6308.5
$ws.Range("H122").Value = 2657.158
$ws.Range("I122").Value = 2660.7
$ws.Range("J122").Value = 2653.2222
$ws.Range("K122").Value = 7982.099999999999
$ws.Range("L122").Value = 7959.6666
$ws.Range("M122").Value = -5532.099999999999
$ws.Range("N122").Value = -12859.6666
$ws.Range("H126").Value = 2416.6667
$ws.Range("I126").Value = 625
$ws.Range("J126").Value = 3312.5
$ws.Range("K126").Value = 1875
$ws.Range("L126").Value = 9937.5
$ws.Range("M126").Value = 595
$ws.Range("N126").Value = -14877.5
$ws.Range("H132").Value = 2548.5117
$ws.Range("I132").Value = 1782.4
$ws.Range("J132").Value = 4316.4614
$ws.Range("K132").Value = 5347.200000000001
$ws.Range("L132").Value = 12949.3842
$ws.Range("M132").Value = -2817.200000000001
$ws.Range("N132").Value = -18009.3842
$ws.Range("H134").Value = 1385.0682
$ws.Range("I134").Value = 921.875
$ws.Range("J134").Value = 1940.9
$ws.Range("K134").Value = 2765.625
$ws.Range("L134").Value = 5822.700000000001
$ws.Range("M134").Value = -230.625
$ws.Range("N134").Value = -10892.7
$ws.Range("H136").Value = 13515759
$ws.Range("I136").Value = 1224
$ws.Range("J136").Value = 41671040
$ws.Range("K136").Value = 3672
$ws.Range("L136").Value = 125013120
$ws.Range("M136").Value = -1122
$ws.Range("N136").Value = -125018220
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 617.8108
$ws.Range("I5").Value = 356.33334
$ws.Range("J5").Value = 2775
$ws.Range("K5").Value = 1069.00002
$ws.Range("L5").Value = 8325
$ws.Range("M5").Value = -957.0000199999999
$ws.Range("N5").Value = -8549
$ws.Range("H92").Value = 1503.5714
$ws.Range("I92").Value = 350
$ws.Range("J92").Value = 2144.4443
$ws.Range("K92").Value = 1050
$ws.Range("L92").Value = 6433.3329
$ws.Range("M92").Value = 198
$ws.Range("N92").Value = -8929.332900000001
$ws.Range("H122").Value = 1392.8889
$ws.Range("I122").Value = 720
$ws.Range("J122").Value = 2065.7778
$ws.Range("K122").Value = 6480
$ws.Range("L122").Value = 18592.0002
$ws.Range("M122").Value = -4030
$ws.Range("N122").Value = -23492.0002
$ws.Range("H135").Value = 617.8108
$ws.Range("I135").Value = 356.33334
$ws.Range("J135").Value = 2775
$ws.Range("K135").Value = 3207.00006
$ws.Range("L135").Value = 24975
$ws.Range("M135").Value = -672.0000600000003
$ws.Range("N135").Value = -30045
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 46929.914
$ws.Range("I68").Value = 30000
$ws.Range("J68").Value = 49469.4
$ws.Range("K68").Value = 30000
$ws.Range("L68").Value = 49469.4
$ws.Range("M68").Value = -29189
$ws.Range("N68").Value = -51091.4
$ws.Range("H71").Value = 46929.914
$ws.Range("I71").Value = 30000
$ws.Range("J71").Value = 49469.4
$ws.Range("K71").Value = 90000
$ws.Range("L71").Value = 148408.2
$ws.Range("M71").Value = -85944
$ws.Range("N71").Value = -156520.2
$ws.Range("H102").Value = 80148.16
$ws.Range("I102").Value = 2592.7
$ws.Range("J102").Value = 338666.34
$ws.Range("K102").Value = 2592.7
$ws.Range("L102").Value = 338666.34
$ws.Range("M102").Value = -970.6999999999998
$ws.Range("N102").Value = -341910.34
$ws.Range("H107").Value = 775.1429000000001
$ws.Range("I107").Value = 315.5
$ws.Range("K107").Value = 315.5
$ws.Range("M107").Value = 1604.5
$ws.Range("H132").Value = 2406.9607
$ws.Range("I132").Value = 1707.3429
$ws.Range("J132").Value = 3937.375
$ws.Range("K132").Value = 5122.028700000001
$ws.Range("L132").Value = 11812.125
$ws.Range("M132").Value = -2592.028700000001
$ws.Range("N132").Value = -16872.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2837.7222
$ws.Range("I132").Value = 1814.1666
$ws.Range("J132").Value = 3861.2778
$ws.Range("K132").Value = 5442.4998
$ws.Range("L132").Value = 11583.8334
$ws.Range("M132").Value = -2912.4998
$ws.Range("N132").Value = -16643.8334
$ws.Range("H136").Value = 2566524
$ws.Range("I136").Value = 4547615
$ws.Range("J136").Value = 2759
$ws.Range("K136").Value = 13642845
$ws.Range("L136").Value = 8277
$ws.Range("M136").Value = -13640295
$ws.Range("N136").Value = -13377
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 180536.97
$ws.Range("I132").Value = 235053
$ws.Range("K132").Value = 705159
$ws.Range("M132").Value = -702629
$ws.Range("H136").Value = 1237.8064
$ws.Range("I136").Value = 698.88
$ws.Range("J136").Value = 3483.3333
$ws.Range("K136").Value = 2096.64
$ws.Range("L136").Value = 10449.9999
$ws.Range("M136").Value = 453.3600000000001
$ws.Range("N136").Value = -15549.9999
